$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 32

$ws.Cells.Item($row, 1).Value = "WGE 341"
$ws.Cells.Item($row, 2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item($row, 3).Value = "16-01-2026"
$ws.Cells.Item($row, 4).Value = 286962
$ws.Cells.Item($row, 5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item($row, 6).Value = 34413429360
$ws.Cells.Item($row, 7).Value = "NEFT"
$ws.Cells.Item($row, 8).Value = "SBIN0003229"
$ws.Cells.Item($row, 9).Value = "AAAFW8862C"
$ws.Cells.Item($row, 10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item($row, 11).Value = "PRADEEP KUMAR K ( ROOM OWNER)"
$ws.Cells.Item($row, 12).Value = "eb2a3bbf-c212-47d8-a926-14b25ccbb4d2"
$ws.Cells.Item($row, 13).Value = 570183114278
$ws.Cells.Item($row, 14).Value = "SBIN 0070200"
# Columns O,P,Q,R,S,T are blank in the source row (kept unset)
$ws.Cells.Item($row, 21).Value = "pending"
$ws.Cells.Item($row, 22).Value = 9000
# Column W is blank in the source row (kept unset)
$ws.Cells.Item($row, 24).Value = "CHELARI ROOM RENT @9000 DUE DATE 16TH OF EVERY MONTH RPA_UNIQUE_ID : 39ebae64-7ca9-4211-97a4-aac615744154"
$ws.Cells.Item($row, 25).Value = "CHELARI PROJECT"
$ws.Cells.Item($row, 26).Value = 0
$ws.Cells.Item($row, 27).Value = "hrm@westernidc.com"
$ws.Cells.Item($row, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item($row, 29).Value = 0
$ws.Cells.Item($row, 30).Value = 0
$ws.Cells.Item($row, 31).Value = 0
# Columns AF..AO are blank in the source row (kept unset)
